$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 21: "7BB-20-6" piezo (2 Pins variant) ---
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = "Piezos"
$ws.Range("D21").Value = "7BB-20-6"
$ws.Range("E21").Value = "2 Pins"
$ws.Range("E21").Style = $ws.Range("E4").Style
$ws.Range("F21").Value = "https://www.distrelec.ch/de/piezo-element-murata-7bb-20/p/13787033?q=*&filter_Category4=Signalgeber+f%C3%BCr+Leiterplatten&filter_Category3=Akustische+Signalgeber&filter_Buyable=1&page=77&origPos=77&origPageSize=25&simi=97.0"
$ws.Range("F21").Style = $ws.Range("F29").Style

# --- Weights (kg) for the piezo rows ---
$ws.Range("E29").Value = 0.378
$ws.Range("E30").Value = 0.2592
$ws.Range("E31").Value = 0.5076
$ws.Range("E32").Value = 1.22
$ws.Range("E33").Value = 0.972
$ws.Range("E34").Value = 0.6912
$ws.Range("E35").Value = 1.3

# --- New hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("F33"), "https://www.distrelec.ch/de/piezo-element-murata-7bb-27/p/13787058?q=7BB-27-4&page=1&origPos=1&origPageSize=25&simi=99.65&no-cache=true")
$ws.Range("F33").Style = $ws.Range("F32").Style

$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.distrelec.ch/de/piezo-element-murata-7bb-20/p/13787033?q=*&filter_Category4=Signalgeber+f%C3%BCr+Leiterplatten&filter_Category3=Akustische+Signalgeber&filter_Buyable=1&page=77&origPos=77&origPageSize=25&simi=97.0")
$ws.Range("F21").Style = $ws.Range("F29").Style

# --- View: scroll back to A1, move selection to E24 ---
$ws.Range("E24").Select()
